# Update CIBMTR Priority Variables 2022 ValueSet metadata sheet:
# - bump version, status, date
# - update contact info (publisher-style contact + named contact)
# - insert a new "Jurisdiction" row before "Description"
# (logo/colors commit message corresponds to this FHIR IG metadata refresh)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version bump
$ws.Range("B3").Value = "0.1.7"

# Status: active -> draft
$ws.Range("B6").Value = "draft"

# Date refresh
$ws.Range("B8").Value = "2024-08-23T10:17:11-05:00"

# Contact details updated
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# Insert a new row for "Jurisdiction" before "Description" (currently row 12)
$ws.Rows("12").Insert()

$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""
$ws.Range("A12:B12").Style = $ws.Range("A13:B13").Style
